# "example has an extra column"
#
# The GeneticsFile example sheet had a stray "OutputGeneValues" column
# (column E) that duplicated information already captured elsewhere. This
# script removes that column - shifting every subsequent column (and the
# header/body values they hold) one slot to the left - then leaves the
# workbook's view state pointed at the GeneticsFile sheet, matching the
# author's final selection there.

$wb = $excel.ActiveWorkbook

$description = $wb.Worksheets.Item("Description")
$genetics = $wb.Worksheets.Item("GeneticsFile")

# Remove the extra "OutputGeneValues" column (E) from the GeneticsFile sheet.
# Everything to the right (OutputNeutralStatistics, OutputFstatsWeirCockerham,
# OutputFstatsWeirHill, OutputStartGenetics, OutputInterval, PatchList,
# NbrPatchesToSample, nIndividualsToSample, Stages) shifts left by one column.
$genetics.Columns("E").Delete() | Out-Null

# Restore the prior selection on the Description sheet (it is no longer the
# active tab once GeneticsFile is selected below).
$description.Activate() | Out-Null
$description.Range("A8").Select() | Out-Null

# GeneticsFile becomes the active sheet/tab, with the selection resting on
# the first data cell of the (now shifted) OutputNeutralStatistics column.
$genetics.Activate() | Out-Null
$genetics.Range("F2").Select() | Out-Null
